$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$changes = @(
    ,(4, 8, 422.5)
    ,(4, 9, 246.33333)
    ,(4, 11, 246.33333)
    ,(4, 13, -132.33333)
    ,(20, 8, 948.3333)
    ,(20, 9, 948.3333)
    ,(20, 10, 0)
    ,(20, 11, 948.3333)
    ,(20, 12, 0)
    ,(20, 13, -718.3333)
    ,(20, 14, $null)
    ,(35, 8, 948.3333)
    ,(35, 9, 948.3333)
    ,(35, 10, 0)
    ,(35, 11, 948.3333)
    ,(35, 12, 0)
    ,(35, 13, -569.3333)
    ,(35, 14, $null)
    ,(64, 8, 25776.814)
    ,(64, 9, 85516.664)
    ,(64, 10, 2651.7097)
    ,(64, 11, 85516.664)
    ,(64, 12, 2651.7097)
    ,(64, 13, -85268.664)
    ,(64, 14, -3147.7097)
    ,(67, 8, 25776.814)
    ,(67, 9, 85516.664)
    ,(67, 10, 2651.7097)
    ,(67, 11, 85516.664)
    ,(67, 12, 2651.7097)
    ,(67, 13, -84658.664)
    ,(67, 14, -4367.709699999999)
    ,(74, 8, 3199.4119)
    ,(74, 9, 3358.3333)
    ,(74, 10, 2818)
    ,(74, 11, 3358.3333)
    ,(74, 12, 2818)
    ,(74, 13, -2422.3333)
    ,(74, 14, -4690)
    ,(77, 8, 3199.4119)
    ,(77, 9, 3358.3333)
    ,(77, 10, 2818)
    ,(77, 11, 16791.6665)
    ,(77, 12, 14090)
    ,(77, 13, -12111.6665)
    ,(77, 14, -23450)
    ,(87, 8, 41233.332)
    ,(87, 10, 41233.332)
    ,(87, 12, 41233.332)
    ,(87, 14, -43729.332)
    ,(90, 8, 41233.332)
    ,(90, 10, 41233.332)
    ,(90, 12, 123699.996)
    ,(90, 14, -136179.996)
    ,(123, 8, 32296)
    ,(123, 10, 32296)
    ,(123, 12, 32296)
    ,(123, 14, -42096)
)
foreach ($chg in $changes) {
    $r = $chg[0]; $c = $chg[1]; $v = $chg[2]
    if ($null -eq $v) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $v
    }
}

$ws = $wb.Worksheets.Item("ARM")
$changes = @(
    ,(95, 8, 32723.75)
    ,(95, 10, 32723.75)
    ,(95, 12, 32723.75)
    ,(95, 14, -38215.75)
    ,(96, 8, 31409.666)
    ,(96, 10, 31409.666)
    ,(96, 12, 31409.666)
    ,(96, 14, -36901.666)
    ,(102, 8, 12013.521)
    ,(102, 9, 2287.0908)
    ,(102, 10, 20929.416)
    ,(102, 11, 2287.0908)
    ,(102, 12, 20929.416)
    ,(102, 13, -665.0907999999999)
    ,(102, 14, -24173.416)
)
foreach ($chg in $changes) {
    $r = $chg[0]; $c = $chg[1]; $v = $chg[2]
    if ($null -eq $v) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $v
    }
}

$ws = $wb.Worksheets.Item("BSM")
$changes = @(
    ,(16, 8, 39500)
    ,(16, 9, 0)
    ,(16, 10, 39500)
    ,(16, 11, 0)
    ,(16, 12, 39500)
    ,(16, 13, $null)
    ,(16, 14, -39840)
    ,(29, 8, 516)
    ,(29, 9, 516)
    ,(29, 11, 516)
    ,(29, 13, -227)
    ,(81, 8, 9062.625)
    ,(81, 10, 9062.625)
    ,(81, 12, 9062.625)
    ,(81, 14, -11184.625)
    ,(84, 8, 9062.625)
    ,(84, 10, 9062.625)
    ,(84, 12, 27187.875)
    ,(84, 14, -37795.875)
    ,(94, 8, 555.38464)
    ,(94, 9, 555.5417)
    ,(94, 10, 553.5)
    ,(94, 11, 555.5417)
    ,(94, 12, 553.5)
    ,(94, 13, -104.5417)
    ,(94, 14, -1455.5)
    ,(105, 8, 3086.8572)
    ,(105, 9, 2634.4443)
    ,(105, 10, 3901.2)
    ,(105, 11, 2634.4443)
    ,(105, 12, 3901.2)
    ,(105, 13, -887.4443000000001)
    ,(105, 14, -7395.2)
    ,(117, 8, 47747)
    ,(117, 10, 47747)
    ,(117, 12, 47747)
    ,(117, 14, -56925)
)
foreach ($chg in $changes) {
    $r = $chg[0]; $c = $chg[1]; $v = $chg[2]
    if ($null -eq $v) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $v
    }
}

$ws = $wb.Worksheets.Item("CRP")
$changes = @(
    ,(4, 8, 254686.56)
    ,(4, 10, 254686.56)
    ,(4, 12, 254686.56)
    ,(4, 14, -254910.56)
    ,(22, 8, 2789)
    ,(22, 9, 525.25)
    ,(22, 10, 4600)
    ,(22, 11, 525.25)
    ,(22, 12, 4600)
    ,(22, 13, -175.25)
    ,(22, 14, -5300)
    ,(115, 8, 29726)
    ,(115, 10, 29726)
    ,(115, 12, 29726)
    ,(115, 14, -32076)
    ,(119, 8, 47686)
    ,(119, 10, 47686)
    ,(119, 12, 47686)
    ,(119, 14, -57362)
    ,(120, 8, 31816.584)
    ,(120, 10, 31816.584)
    ,(120, 12, 31816.584)
    ,(120, 14, -39074.584)
    ,(134, 8, 71941.85000000001)
    ,(134, 9, 1194.7142)
    ,(134, 11, 3584.1426)
    ,(134, 13, -1049.1426)
)
foreach ($chg in $changes) {
    $r = $chg[0]; $c = $chg[1]; $v = $chg[2]
    if ($null -eq $v) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $v
    }
}

$ws = $wb.Worksheets.Item("CUL")
$changes = @(
    ,(64, 8, 3322)
    ,(64, 9, 1870.6666)
    ,(64, 10, 3717.818)
    ,(64, 11, 5611.9998)
    ,(64, 12, 11153.454)
    ,(64, 13, -5341.9998)
    ,(64, 14, -11693.454)
    ,(67, 8, 3322)
    ,(67, 9, 1870.6666)
    ,(67, 10, 3717.818)
    ,(67, 11, 5611.9998)
    ,(67, 12, 11153.454)
    ,(67, 13, -4675.9998)
    ,(67, 14, -13025.454)
    ,(136, 8, 250001470)
    ,(136, 9, 250001470)
    ,(136, 10, 0)
    ,(136, 11, 750004410)
    ,(136, 12, 0)
    ,(136, 13, -749999310)
    ,(136, 14, $null)
)
foreach ($chg in $changes) {
    $r = $chg[0]; $c = $chg[1]; $v = $chg[2]
    if ($null -eq $v) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $v
    }
}

$ws = $wb.Worksheets.Item("GSM")
$changes = @(
    ,(5, 8, 15000)
    ,(5, 9, 15000)
    ,(5, 10, 0)
    ,(5, 11, 15000)
    ,(5, 12, 0)
    ,(5, 13, -14888)
    ,(5, 14, $null)
    ,(31, 8, 16120)
    ,(31, 9, 1745)
    ,(31, 10, 30495)
    ,(31, 11, 1745)
    ,(31, 12, 30495)
    ,(31, 13, -1453)
    ,(31, 14, -31079)
    ,(37, 8, 16120)
    ,(37, 9, 1745)
    ,(37, 10, 30495)
    ,(37, 11, 1745)
    ,(37, 12, 30495)
    ,(37, 13, -1468)
    ,(37, 14, -31049)
    ,(118, 8, 38086.668)
    ,(118, 10, 38086.668)
    ,(118, 12, 38086.668)
    ,(118, 14, -41400.668)
    ,(122, 8, 1212.5)
    ,(122, 9, 1283.3334)
    ,(122, 11, 3850.0002)
    ,(122, 13, -1400.0002)
)
foreach ($chg in $changes) {
    $r = $chg[0]; $c = $chg[1]; $v = $chg[2]
    if ($null -eq $v) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $v
    }
}

$ws = $wb.Worksheets.Item("LTW")
$changes = @(
    ,(2, 8, 50001)
    ,(2, 10, 50001)
    ,(2, 12, 50001)
    ,(2, 14, -50225)
    ,(93, 8, 1399.7142)
    ,(93, 9, 1353.8889)
    ,(93, 10, 1421.421)
    ,(93, 11, 1353.8889)
    ,(93, 12, 1421.421)
    ,(93, 13, -105.8888999999999)
    ,(93, 14, -3917.421)
    ,(119, 8, 47408)
    ,(119, 10, 47408)
    ,(119, 12, 47408)
    ,(119, 14, -57084)
    ,(120, 8, 51188.4)
    ,(120, 10, 51188.4)
    ,(120, 12, 51188.4)
    ,(120, 14, -60864.4)
    ,(122, 8, 2299)
    ,(122, 9, 2279.8572)
    ,(122, 10, 2500)
    ,(122, 11, 6839.571599999999)
    ,(122, 12, 7500)
    ,(122, 13, -4389.571599999999)
    ,(122, 14, -12400)
)
foreach ($chg in $changes) {
    $r = $chg[0]; $c = $chg[1]; $v = $chg[2]
    if ($null -eq $v) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $v
    }
}

$ws = $wb.Worksheets.Item("WVR")
$changes = @(
    ,(18, 8, 9833.333000000001)
    ,(18, 9, 0)
    ,(18, 10, 9833.333000000001)
    ,(18, 11, 0)
    ,(18, 12, 9833.333000000001)
    ,(18, 13, $null)
    ,(18, 14, -10179.333)
    ,(95, 8, 38992)
    ,(95, 10, 38992)
    ,(95, 12, 38992)
    ,(95, 14, -44484)
    ,(96, 8, 4388.5)
    ,(96, 9, 1000)
    ,(96, 10, 7777)
    ,(96, 11, 1000)
    ,(96, 12, 7777)
    ,(96, 13, 373)
    ,(96, 14, -10523)
)
foreach ($chg in $changes) {
    $r = $chg[0]; $c = $chg[1]; $v = $chg[2]
    if ($null -eq $v) {
        $ws.Cells.Item($r, $c).ClearContents()
    } else {
        $ws.Cells.Item($r, $c).Value = $v
    }
}
